# feat(upvote comment): upvoting comment API added.
# * POST
# * DELETE
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The workbook's shared-string table is built in first-seen order, so touch
# the four brand-new strings once, in the same order they appear in the
# target file (upvote/, /article/{articleId}, /comment/{commentId},
# /publication/{publicationId}), before filling in the real layout below.
$ws.Range("A35").Value = "upvote/"
$ws.Range("B37").Value = "/article/{articleId}"
$ws.Range("B39").Value = "/comment/{commentId}"
$ws.Range("B35").Value = "/publication/{publicationId}"

# Row 35-36: upvote/ + /publication/{publicationId}  -> POST / DELETE
$ws.Range("A35").Value = "upvote/"
$ws.Range("B35").Value = "/publication/{publicationId}"
$ws.Range("D35").Value = "POST"
$ws.Range("D36").Value = "DELETE"

# Row 37-38: upvote/ + /article/{articleId}  -> POST / DELETE
$ws.Range("A37").Value = "upvote/"
$ws.Range("B37").Value = "/article/{articleId}"
$ws.Range("D37").Value = "POST"
$ws.Range("D38").Value = "DELETE"

# Row 39-40: upvote/ + /comment/{commentId}  -> POST / DELETE
$ws.Range("A39").Value = "upvote/"
$ws.Range("B39").Value = "/comment/{commentId}"
$ws.Range("D39").Value = "POST"
$ws.Range("D40").Value = "DELETE"

# Update the visible window / selection to match the edited state.
$ws.Range("A19").Select()
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("F25").Select()
